$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "dvno" column (D) values for rows 2-6
$ws.Range("D2").Value = 7080
$ws.Range("D3").Value = 7021
$ws.Range("D4").Value = 7022
$ws.Range("D5").Value = 7023
$ws.Range("D6").Value = 7024

# Update the active selection to D6
$ws.Activate()
$ws.Range("D6").Select()
